# Updated with new color settings
# Converts Subset labels in column A from the legacy
# PFIZER/PFIZER_BIVALENT/MODERNA/MODERNA_BIVALENT/JANSSEN/NOVAVAX token
# format to the new Pfizer_mono/Pfizer_bi/Moderna_mono/Moderna_bi/Janssen/Novavax
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Convert-SubsetName {
    param([string]$s)

    $tokens = $s -split "_"
    $result = @()
    $i = 0
    while ($i -lt $tokens.Count) {
        $tok = $tokens[$i]
        if ($tok -eq "PFIZER" -or $tok -eq "MODERNA") {
            $name = $tok.Substring(0,1) + $tok.Substring(1).ToLower()
            if (($i + 1) -lt $tokens.Count -and $tokens[$i + 1] -eq "BIVALENT") {
                $result += "${name}_bi"
                $i += 2
            } else {
                $result += "${name}_mono"
                $i += 1
            }
        } elseif ($tok -eq "JANSSEN") {
            $result += "Janssen"
            $i += 1
        } elseif ($tok -eq "NOVAVAX") {
            $result += "Novavax"
            $i += 1
        } else {
            $result += $tok
            $i += 1
        }
    }
    return ($result -join "_")
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# NOTE: this engine's -eq/-ne/-ceq/-cne string comparison operators are
# case-insensitive, which would incorrectly skip rows such as
# "JANSSEN_NOVAVAX" -> "Janssen_Novavax" (values differing only by case)
# if used to decide whether to write back. So we unconditionally
# reassign every non-blank cell instead of conditioning the write on a
# changed/unchanged comparison.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null -and $old -ne "") {
        $new = Convert-SubsetName $old
        $cell.Value = $new
    }
}
